$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet/tab: date moves from 2022-05-22 to 2022-05-23
$ws.Name = "Through 2022-05-23"

# Update header label in I1 (shared string) to reflect new "through" date
$ws.Range("I1").Value = "2022 (through 05-23)"

# Update newly reported data for May 2022 (row 6) and the Total row (row 14)
$ws.Range("I6").Value = 83
$ws.Range("I14").Value = 635
